# Generate Report for Handoff
# Replaces the old source-file UUID (6c82ee76-...) with a new one
# (dbb4ade7-...), refreshes the handoff/handback timestamps + xliff file
# names, and clears out the "Latest Target File" / "Latest Handback File"
# columns on the per-locale sheets (the handback step hasn't produced a
# target/handback file yet for the new source file).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
$ov.Range("B2").Value = "e2e\dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
$ov.Range("G2").Value = "2016-08-28 00:57:19"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
$zh.Range("G2").Value = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.038d6863be4d02522a5fc2b16c830271a3cdeb05.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-28 00:57:14"
$zh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$zh.Range("I2").ClearFormats()
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""

$zh.Columns.Item(9).ColumnWidth = 17.8
$zh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
$de.Range("G2").Value = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.038d6863be4d02522a5fc2b16c830271a3cdeb05.de-de.xlf"
$de.Range("H2").Value = "2016-08-28 00:57:19"
$de.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "dbb4ade7-66c5-4e1c-8025-87d5578914e7.md"
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$de.Range("I2").ClearFormats()
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""

$de.Columns.Item(9).ColumnWidth = 17.8
$de.Columns.Item(10).ColumnWidth = 20.8
